$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.463
$ws.Range("A3").Value = -21.516
$ws.Range("C5").Value = -12.836
$ws.Range("E5").Value = 13.034
$ws.Range("E9").Value = 12.82
$ws.Range("E11").Value = 13.068
$ws.Range("A14").Value = -20.828
$ws.Range("A21").Value = -21.04
$ws.Range("E21").Value = 13.535
$ws.Range("A23").Value = -21.709
$ws.Range("A25").Value = -22.269
